$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision process has concluded without reaching an agreement on a movie for Friday.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as no consensus was reached regarding the movie to show on Friday.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision about the movie for Friday has not been made.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected during the meeting.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" has been successfully recorded.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding the movie to show on Friday.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday ended without a definitive choice being made.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision for Friday's movie was not finalized, so we have no selection to make.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to be shown on Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie will be shown on Friday based on the committee's discussion.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The rights to both movies, `"Barbie`" and `"Oppenheimer,`" have been successfully acquired.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision about what movie to play on Friday resulted in no agreement.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision process concluded without a consensus on which movie to show on Friday, so I have recorded the outcome as no decision made.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision,`" indicating that there was no agreement on which movie to show on Friday.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("D18").Value = "Oppenheimer_was_selected, "
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision not to acquire a movie for Friday has been recorded.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday was not made by the committee.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for Friday's showing.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding which movie to show on Friday.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has been recorded as `"no_decision.`"`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("D25").Value = "Barbie_was_selected, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday was not reached, and thus the proper course of action was to conclude that there is no decision.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be shown on Friday.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been successfully recorded with `"Barbie`" being selected as the movie to acquire for the event.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision-making process ended without a clear choice for Friday's movie, therefore no action was taken regarding the acquisition of movie rights.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie selected for Friday is `"Barbie.`"`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been selected for the assembly.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights to `"Oppenheimer`" will be acquired for the upcoming screening.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not concluded, so I am calling the no_decision function.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to show `"Barbie`" on Friday.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday was concluded without an agreement, resulting in no decision being made.`n"
